$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 2063.44
$ws.Range("J112").Value = 2044.9
$ws.Range("L112").Value = 6134.700000000001
$ws.Range("N112").Value = -8350.700000000001
$ws.Range("H121").Value = 2996
$ws.Range("J121").Value = 2996
$ws.Range("L121").Value = 8988
$ws.Range("N121").Value = -12482
$ws.Range("H132").Value = 102518.5
$ws.Range("I132").Value = 127147.875
$ws.Range("J132").Value = 4001
$ws.Range("K132").Value = 381443.625
$ws.Range("L132").Value = 12003
$ws.Range("M132").Value = -378913.625
$ws.Range("N132").Value = -17063
$ws.Range("H137").Value = 2164.2632
$ws.Range("I137").Value = 2191.5334
$ws.Range("J137").Value = 2062
$ws.Range("K137").Value = 6574.600199999999
$ws.Range("L137").Value = 6186
$ws.Range("M137").Value = -4024.600199999999
$ws.Range("N137").Value = -11286
$ws.Range("H138").Value = 2375.3484
$ws.Range("I138").Value = 981.2381
$ws.Range("K138").Value = 2943.7143
$ws.Range("M138").Value = 2196.2857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 13219.8
$ws.Range("I28").Value = 13219.8
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 13219.8
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -13027.8
$ws.Range("N28").ClearContents()
$ws.Range("H32").Value = 1414.4117
$ws.Range("I32").Value = 1414.4117
$ws.Range("K32").Value = 1414.4117
$ws.Range("M32").Value = -1127.4117
$ws.Range("H61").Value = 3891.182
$ws.Range("I61").Value = 3891.182
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3891.182
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3679.182
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 2014532.6
$ws.Range("I74").Value = 1030471.8
$ws.Range("K74").Value = 1030471.8
$ws.Range("M74").Value = -1029597.8
$ws.Range("H77").Value = 2014532.6
$ws.Range("I77").Value = 1030471.8
$ws.Range("K77").Value = 5152359
$ws.Range("M77").Value = -5147991
$ws.Range("H99").Value = 13219.8
$ws.Range("I99").Value = 13219.8
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 13219.8
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -10224.8
$ws.Range("N99").ClearContents()
$ws.Range("H102").Value = 556.5714
$ws.Range("I102").Value = 556.5714
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 556.5714
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 1065.4286
$ws.Range("N102").ClearContents()
$ws.Range("H132").Value = 15625864
$ws.Range("I132").Value = 891.6129
$ws.Range("K132").Value = 2674.8387
$ws.Range("M132").Value = -144.8386999999998
$ws.Range("H136").Value = 3891.182
$ws.Range("I136").Value = 3891.182
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 11673.546
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -9123.545999999998
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 7498.5
$ws.Range("I26").Value = 7498.5
$ws.Range("K26").Value = 7498.5
$ws.Range("M26").Value = -7206.5
$ws.Range("H35").Value = 24582
$ws.Range("J35").Value = 38592
$ws.Range("L35").Value = 38592
$ws.Range("N35").Value = -39212
$ws.Range("H94").Value = 1460.4375
$ws.Range("I94").Value = 1585.5
$ws.Range("K94").Value = 1585.5
$ws.Range("M94").Value = -1134.5
$ws.Range("H134").Value = 44873444
$ws.Range("I134").Value = 20835118
$ws.Range("J134").Value = 333333340
$ws.Range("K134").Value = 62505354
$ws.Range("L134").Value = 1000000020
$ws.Range("M134").Value = -62502819
$ws.Range("N134").Value = -1000005090

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 829.6
$ws.Range("I22").Value = 799.5
$ws.Range("J22").Value = 849.6667
$ws.Range("K22").Value = 799.5
$ws.Range("L22").Value = 849.6667
$ws.Range("M22").Value = -449.5
$ws.Range("N22").Value = -1549.6667
$ws.Range("H31").Value = 2556.2964
$ws.Range("I31").Value = 2015.7858
$ws.Range("K31").Value = 2015.7858
$ws.Range("M31").Value = -1720.7858
$ws.Range("H34").Value = 2556.2964
$ws.Range("I34").Value = 2015.7858
$ws.Range("K34").Value = 2015.7858
$ws.Range("M34").Value = -1813.7858
$ws.Range("H132").Value = 4615.1035
$ws.Range("I132").Value = 3670.423
$ws.Range("J132").Value = 12802.333
$ws.Range("K132").Value = 11011.269
$ws.Range("L132").Value = 38406.999
$ws.Range("M132").Value = -8481.269
$ws.Range("N132").Value = -43466.999
$ws.Range("H134").Value = 3032458.8
$ws.Range("I134").Value = 1936.0769
$ws.Range("J134").Value = 14288686
$ws.Range("K134").Value = 5808.2307
$ws.Range("L134").Value = 42866058
$ws.Range("M134").Value = -3273.2307
$ws.Range("N134").Value = -42871128

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 148686
$ws.Range("I46").Value = 500451
$ws.Range("K46").Value = 1501353
$ws.Range("M46").Value = -1501262
$ws.Range("H113").Value = 582.06665
$ws.Range("I113").Value = 387.5
$ws.Range("J113").Value = 652.8182
$ws.Range("K113").Value = 1162.5
$ws.Range("L113").Value = 1958.4546
$ws.Range("M113").Value = 1007.5
$ws.Range("N113").Value = -6298.4546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9525.784
$ws.Range("J70").Value = 9742.375
$ws.Range("L70").Value = 9742.375
$ws.Range("N70").Value = -10282.375
$ws.Range("H73").Value = 9525.784
$ws.Range("J73").Value = 9742.375
$ws.Range("L73").Value = 9742.375
$ws.Range("N73").Value = -11614.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2371.1538
$ws.Range("I16").Value = 2365.75
$ws.Range("J16").Value = 2379.8
$ws.Range("K16").Value = 2365.75
$ws.Range("L16").Value = 2379.8
$ws.Range("M16").Value = -2195.75
$ws.Range("N16").Value = -2719.8
$ws.Range("H21").Value = 14333.333
$ws.Range("J21").Value = 14333.333
$ws.Range("L21").Value = 14333.333
$ws.Range("N21").Value = -14681.333
$ws.Range("H22").Value = 1243.25
$ws.Range("I22").Value = 1079.3125
$ws.Range("J22").Value = 1899
$ws.Range("K22").Value = 1079.3125
$ws.Range("L22").Value = 1899
$ws.Range("M22").Value = -784.3125
$ws.Range("N22").Value = -2489
$ws.Range("H27").Value = 1243.25
$ws.Range("I27").Value = 1079.3125
$ws.Range("J27").Value = 1899
$ws.Range("K27").Value = 1079.3125
$ws.Range("L27").Value = 1899
$ws.Range("M27").Value = -972.3125
$ws.Range("N27").Value = -2113
$ws.Range("H55").Value = 419.27777
$ws.Range("I55").Value = 315.4375
$ws.Range("K55").Value = 315.4375
$ws.Range("M55").Value = -142.4375
$ws.Range("H61").Value = 2675.5
$ws.Range("I61").Value = 1780.8
$ws.Range("J61").Value = 4166.6665
$ws.Range("K61").Value = 1780.8
$ws.Range("L61").Value = 4166.6665
$ws.Range("M61").Value = -1578.8
$ws.Range("N61").Value = -4570.6665
$ws.Range("H100").Value = 607599.2
$ws.Range("I100").Value = 607599.2
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 607599.2
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -607058.2
$ws.Range("N100").ClearContents()
$ws.Range("H113").Value = 2675.5
$ws.Range("I113").Value = 1780.8
$ws.Range("J113").Value = 4166.6665
$ws.Range("K113").Value = 1780.8
$ws.Range("L113").Value = 4166.6665
$ws.Range("M113").Value = 389.2
$ws.Range("N113").Value = -8506.666499999999
$ws.Range("H125").Value = 69666
$ws.Range("J125").Value = 69499
$ws.Range("L125").Value = 69499
$ws.Range("N125").Value = -79339
$ws.Range("H136").Value = 20002720
$ws.Range("I136").Value = 2612.8293
$ws.Range("J136").Value = 111114320
$ws.Range("K136").Value = 7838.4879
$ws.Range("L136").Value = 333342960
$ws.Range("M136").Value = -5288.4879
$ws.Range("N136").Value = -333348060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 7638.3335
$ws.Range("I51").Value = 5468.25
$ws.Range("K51").Value = 5468.25
$ws.Range("M51").Value = -4958.25
$ws.Range("H132").Value = 2123.375
$ws.Range("I132").Value = 1997.5
$ws.Range("J132").Value = 2249.25
$ws.Range("K132").Value = 5992.5
$ws.Range("L132").Value = 6747.75
$ws.Range("M132").Value = -3462.5
$ws.Range("N132").Value = -11807.75
$ws.Range("H136").Value = 656.9259
$ws.Range("I136").Value = 586.1923
$ws.Range("J136").Value = 2496
$ws.Range("K136").Value = 1758.5769
$ws.Range("L136").Value = 7488
$ws.Range("M136").Value = 791.4231
$ws.Range("N136").Value = -12588
